# Fill in the two pending TxHash placeholder sheets ("B1" and "B2") with
# the real transfer hashes, and move the active tab from "A12" over to "B2"
# (mirrors a user finishing data entry on the last-worked-on sheet).
#
# NOTE on write order: the engine appends brand-new shared strings to the
# table in first-write order, so cells are touched in the exact sequence
# that reproduces the target shared-string layout (B1!A3, B1!A2, then
# B2!A2, B2!A3).

$wb = $excel.ActiveWorkbook

# --- "B1" sheet: fill in the two TxHash cells ---
$ws1 = $wb.Worksheets.Item("B1")
$ws1.Range("A3").Value = "8D0224B36F0C6BC16C094E705BC9849A799741EA666174482286458A801F92DE"
$ws1.Range("A2").Value = "D0300EDE06E8B58E80447493E4C463D9FB2150A92613A23861FFF93EBD27B037"
$ws1.Activate()
$ws1.Range("A2").Select()

# --- "B2" sheet: fill in the two TxHash cells ---
$ws2 = $wb.Worksheets.Item("B2")
$ws2.Range("A2").Value = "87F4590CE9B978449BF027DB44994F18B047D6D28FBAB99D5DF8A881AF3EA71E"
$ws2.Range("A3").Value = "4DB60BB72A1D87B5A84487488BC8742FFF95807183D18904A8FF57739EDF735B"

# "B2" ends up the active/selected sheet & cell, as in the source workbook.
$ws2.Activate()
$ws2.Range("A3").Select()
